# "Generate Report for handback"
#
# Refreshes the handoff/handback timestamps recorded on the per-locale
# report sheets (mirrors a re-run of the handback report generator).
#
#   zh-cn sheet: D2 (Correspond Handoff Datetime)  -> 2016-01-13 04:48:14
#                G2 (Correspond Handback DateTime)  -> 2016-01-13 04:49:20
#   de-de sheet: D2 (Correspond Handoff Datetime)  -> 2016-01-13 04:48:35
#                G2 (Correspond Handback DateTime)  -> 2016-01-13 04:49:53

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-13 04:48:14"
$wsZhCn.Range("G2").Value = "2016-01-13 04:49:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-13 04:48:35"
$wsDeDe.Range("G2").Value = "2016-01-13 04:49:53"
